$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 6); remaining rows keep their own data for now
$ws.Rows.Item(6).Delete()

# Adjust custom column widths (col I: 7 -> 8, col Q: 8 -> 7) expressed in
# Excel "ColumnWidth" character units; COM pads by ~0.8333 vs the raw OOXML width.
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 6.166666666666667

# Replace the 4 data rows (rows 2-5) with the new dataset
# Row 2
$ws.Range("A2").Value = 45140.50694444445
$ws.Range("B2").Value = 4.284
$ws.Range("C2").Value = 5.265
$ws.Range("D2").Value = 0.768
$ws.Range("E2").Value = 10.407
$ws.Range("F2").Value = 9.715
$ws.Range("G2").Value = 5.079
$ws.Range("H2").Value = 10.361
$ws.Range("I2").Value = 6.858
$ws.Range("J2").Value = 3.662
$ws.Range("K2").Value = 6.216
$ws.Range("L2").Value = 7.246
$ws.Range("M2").Value = 4.306
$ws.Range("N2").Value = 0.973
$ws.Range("O2").Value = 3.629
$ws.Range("P2").Value = 6.538
$ws.Range("Q2").Value = 2.613
$ws.Range("R2").Value = 0.117
$ws.Range("S2").Value = 0.053
$ws.Range("T2").Value = 60.452
$ws.Range("U2").Value = 12.195
$ws.Range("V2").Value = 3.56
$ws.Range("W2").Value = 7.364
$ws.Range("X2").Value = 6.727
$ws.Range("Y2").Value = 0.825
$ws.Range("Z2").Value = 5.701
$ws.Range("AA2").Value = 3.6
$ws.Range("AB2").Value = 5.538
$ws.Range("AC2").Value = 2.243
$ws.Range("AD2").Value = 6.732
$ws.Range("AE2").Value = 0.372
$ws.Range("AF2").Value = 8.518000000000001
$ws.Range("AG2").Value = 2.17
$ws.Range("AH2").Value = 5.204

# Row 3
$ws.Range("A3").Value = 45140.51388888889
$ws.Range("B3").Value = 15.051
$ws.Range("C3").Value = 12.305
$ws.Range("D3").Value = 0.852
$ws.Range("E3").Value = 33.459
$ws.Range("F3").Value = 28.184
$ws.Range("G3").Value = 12.651
$ws.Range("H3").Value = 42.98
$ws.Range("I3").Value = 19.275
$ws.Range("J3").Value = 8.927
$ws.Range("K3").Value = 13.611
$ws.Range("L3").Value = 15.022
$ws.Range("M3").Value = 14.198
$ws.Range("N3").Value = 3.711
$ws.Range("O3").Value = 11.876
$ws.Range("P3").Value = 17.842
$ws.Range("Q3").Value = 9.853999999999999
$ws.Range("R3").Value = 0.148
$ws.Range("S3").Value = 0.32
$ws.Range("T3").Value = 181.145
$ws.Range("U3").Value = 34.59
$ws.Range("V3").Value = 11.049
$ws.Range("W3").Value = 23.079
$ws.Range("X3").Value = 13.398
$ws.Range("Y3").Value = 1.704
$ws.Range("Z3").Value = 22.168
$ws.Range("AA3").Value = 10.088
$ws.Range("AB3").Value = 9.935
$ws.Range("AC3").Value = 9.726000000000001
$ws.Range("AD3").Value = 15.245
$ws.Range("AE3").Value = 0.237
$ws.Range("AF3").Value = 38.956
$ws.Range("AG3").Value = 6.474
$ws.Range("AH3").Value = 14.317

# Row 4
$ws.Range("A4").Value = 45140.52083333334
$ws.Range("B4").Value = 4.242
$ws.Range("C4").Value = 3.857
$ws.Range("D4").Value = 0.379
$ws.Range("E4").Value = 9.773
$ws.Range("F4").Value = 8.42
$ws.Range("G4").Value = 3.868
$ws.Range("H4").Value = 18.655
$ws.Range("I4").Value = 5.911
$ws.Range("J4").Value = 2.876
$ws.Range("K4").Value = 4.377
$ws.Range("L4").Value = 5.007
$ws.Range("M4").Value = 4.161
$ws.Range("N4").Value = 1.023
$ws.Range("O4").Value = 3.381
$ws.Range("P4").Value = 5.524
$ws.Range("Q4").Value = 2.885
$ws.Range("R4").Value = 0.046
$ws.Range("S4").Value = 0.012
$ws.Range("T4").Value = 49.711
$ws.Range("U4").Value = 10.767
$ws.Range("V4").Value = 3.165
$ws.Range("W4").Value = 7.113
$ws.Range("X4").Value = 4.475
$ws.Range("Y4").Value = 0.531
$ws.Range("Z4").Value = 9.061999999999999
$ws.Range("AA4").Value = 3.051
$ws.Range("AB4").Value = 3.382
$ws.Range("AC4").Value = 2.699
$ws.Range("AD4").Value = 4.865
$ws.Range("AE4").Value = 0.175
$ws.Range("AF4").Value = 17.378
$ws.Range("AG4").Value = 1.953
$ws.Range("AH4").Value = 4.341

# Row 5
$ws.Range("A5").Value = 45140.52777777778
$ws.Range("B5").Value = 2.45
$ws.Range("C5").Value = 2.35
$ws.Range("D5").Value = 0.26
$ws.Range("E5").Value = 5.77
$ws.Range("F5").Value = 5.02
$ws.Range("G5").Value = 2.32
$ws.Range("H5").Value = 11.02
$ws.Range("I5").Value = 3.58
$ws.Range("J5").Value = 1.75
$ws.Range("K5").Value = 2.69
$ws.Range("L5").Value = 3.14
$ws.Range("M5").Value = 2.46
$ws.Range("N5").Value = 0.58
$ws.Range("O5").Value = 1.96
$ws.Range("P5").Value = 3.31
$ws.Range("Q5").Value = 1.72
$ws.Range("R5").Value = 0.02
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 26.85
$ws.Range("U5").Value = 6.42
$ws.Range("V5").Value = 1.84
$ws.Range("W5").Value = 4.21
$ws.Range("X5").Value = 2.78
$ws.Range("Y5").Value = 0.32
$ws.Range("Z5").Value = 5.28
$ws.Range("AA5").Value = 1.83
$ws.Range("AB5").Value = 2.13
$ws.Range("AC5").Value = 1.55
$ws.Range("AD5").Value = 2.98
$ws.Range("AE5").Value = 0.14
$ws.Range("AF5").Value = 10.13
$ws.Range("AG5").Value = 1.17
$ws.Range("AH5").Value = 2.6
